$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the new data row (row 2) -----------------------------------------
$ws.Range("A2").Value = "MCH343"
$ws.Range("C2").Value = "CASES & AFFIDAVITS, EXTRACTS FROM TRIAL RECORDS OF TONY WEAVERS TRIAL (TO BE SORTED"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 33I | GRAP COUNT NUMER: NONE"

# Cells D2 and H2 stay empty but still carry the row's formatting, just like
# the other record cells.
$dataCells = @("A2", "C2", "D2", "E2", "F2", "G2", "H2")
foreach ($addr in $dataCells) {
    $rng = $ws.Range($addr)
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 10
    $rng.Font.ThemeColor = 1
}

# F2 (extentAndMedium) ends up with its own (alignment-flagged) style.
$ws.Range("F2").WrapText = $false

# --- Restore the frozen header pane and selection --------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("G17").Select()
